# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets
# to reflect the latest scrape, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Map of sheet name -> cell address -> new value
$updates = @{
    "展览" = @{
        "F2"  = 13804
        "F3"  = 90
        "F7"  = 1207
        "F9"  = 13895
        "F10" = 14765
        "F12" = 5
        "F20" = 19
        "F26" = 5732
        "F29" = 5415
        "F32" = 263
    }
    "全部类型" = @{
        "F2"  = 13804
        "F3"  = 90
        "F8"  = 1207
        "F10" = 13895
        "F11" = 14765
        "F13" = 5
        "F21" = 19
        "F27" = 5732
        "F30" = 5415
        "F33" = 263
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($addr in $cellUpdates.Keys) {
        $ws.Range($addr).Value = $cellUpdates[$addr]
    }
}
